$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A68").Value = 43970
$ws.Range("B68").Value = 535
$ws.Range("C68").Value = 157
$ws.Range("D68").Value = 291
$ws.Range("E68").Value = 13
$ws.Range("F68").Value = 17

$ws.Range("A67:F67").Copy()
$ws.Range("A68:F68").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F68"))

$ws.Range("F69").Select()
